$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 3
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 12
$ws.Range("H3").Value = 15

# Row 9
$ws.Range("E9").Value = 18
$ws.Range("F9").Value = 7
$ws.Range("H9").Value = 11

# Row 10
$ws.Range("E10").Value = 31

# Row 11
$ws.Range("E11").Value = 18

# Row 15
$ws.Range("E15").Value = 97
$ws.Range("F15").Value = 46
$ws.Range("H15").Value = 57

# Row 16
$ws.Range("E16").Value = 313

# Row 18
$ws.Range("E18").Value = 96
